# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages update).

$wb = $excel.ActiveWorkbook

# Map of row number -> new F-column value, identical on both sheets.
$updates = @{
    3  = 8347
    5  = 167
    6  = 212
    8  = 761
    10 = 5505
    11 = 10
    17 = 167
    18 = 219
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}

$wb.Save()
